$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.829.12'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.661.38'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.592'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.10%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.65'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.356'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.60'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.133.62'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.718.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.648.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.42'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '343.65'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.79'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.63'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +13.76%  '
$ws.Range('E25').Value = '  +5.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '575.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +23.01%  '
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('E28').Value = '  +2.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.08'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('E31').Value = '  +3.24%  '
$ws.Range('E32').Value = '  +11.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0₃0820'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '175.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.402'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.73'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.19'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.76'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '171.78'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '40.46'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.77'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.631'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0553'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0240'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.30%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0962'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0217'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +12.14%  '
